$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 45

$ws.Range("A45:D45").NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-14"
$ws.Cells.Item($row, 2).Value = "10:55:34"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "24"

$ws.Cells.Item($row, 5).Value = 121136
$ws.Cells.Item($row, 6).Value = 135140
$ws.Cells.Item($row, 7).Value = 161079
$ws.Cells.Item($row, 8).Value = 132400
$ws.Cells.Item($row, 9).Value = 176482
$ws.Cells.Item($row, 10).Value = 113872
$ws.Cells.Item($row, 11).Value = 202268
$ws.Cells.Item($row, 12).Value = 223182
$ws.Cells.Item($row, 13).Value = 173839
$ws.Cells.Item($row, 14).Value = 101866
$ws.Cells.Item($row, 15).Value = 38873
$ws.Cells.Item($row, 16).Value = 34030
$ws.Cells.Item($row, 17).Value = 51481
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36447
$ws.Cells.Item($row, 20).Value = -1

$ws.Range("A45:D45").ClearFormats()
